$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 86

$ws.Cells.Item($row, 1).Value = "2024-10-27 00:00:00"
$ws.Cells.Item($row, 2).Value = 73300
$ws.Cells.Item($row, 3).Value = 10273.44
$ws.Cells.Item($row, 4).Value = 9091.540000000001
$ws.Cells.Item($row, 5).Value = 7.1211
